$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 0.0000000000000000000000001009048904655203
$ws.Range("E2").Value = 0.0000000000000000000000001009048904655203

# Row 3
$ws.Range("C3").Value = $false
$ws.Range("D3").Value = 0.9880664884257978
$ws.Range("E3").Value = 0.9880664884257978

# Row 4
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 0

# Row 5
$ws.Range("D5").Value = 0.004141732160331217
$ws.Range("E5").Value = 0.9958582678396688

# Row 6
$ws.Range("D6").Value = 0.9999886441459673
$ws.Range("E6").Value = 0.00001135585403266859

# Row 7
$ws.Range("D7").Value = 0.00000000000000001095249105107897
$ws.Range("E7").Value = 1

# Row 8
$ws.Range("D8").Value = 0.00000000000000000000000000000000000000000000002799138329023281
$ws.Range("F8").Value = 22.30860328674316
$ws.Range("G8").Value = 0.4285714285714285
